$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updated values, keyed by row number
$values = @{
    2  = 1
    3  = -1
    4  = -8
    5  = -3
    6  = -3
    7  = 1
    8  = -5
    9  = -3
    11 = -1
    12 = 5
    13 = 4
    15 = -2
    16 = 2
    17 = 3
    18 = -3
    19 = -1
    20 = -2
    21 = -2
    25 = 6
    26 = -5
    27 = 0
    28 = 2
    29 = 10
    30 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
